$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-1651168713429141"
$ws1.Range("B2").Value = "go_stims-16511687133874135.csv"
$ws1.Range("B3").Value = "GNG_stims-16511687134102376.csv"
$ws1.Range("B4").Value = "go_stims-16511687134123068.csv"
$ws1.Range("B5").Value = "GNG_stims-16511687134272408.csv"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16511687157842805"
$ws2.Range("B2").Value = "ZB-match_1-16511687140644267.csv"
$ws2.Range("B3").Value = "ZB-match_3-1651168714017146.csv"
$ws2.Range("B4").Value = "TB-16511687154333007.csv"
$ws2.Range("B5").Value = "TB-16511687148053193.csv"
$ws2.Range("B6").Value = "OB-1651168714242898.csv"
$ws2.Range("B7").Value = "ZB-match_8-16511687135532124.csv"
$ws2.Range("B8").Value = "OB-16511687145394976.csv"
$ws2.Range("B9").Value = "OB-16511687142113855.csv"
$ws2.Range("B10").Value = "TB-16511687157588246.csv"

$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-1651168715785246"
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16511687158315954"
$ws4.Range("B2").Value = "MM_stims-16511687157998545.csv"
$ws4.Range("B3").Value = "ZM_stims-16511687157862463.csv"
$ws4.Range("B4").Value = "MM_stims-16511687158154333.csv"
$ws4.Range("B5").Value = "ZM_stims-16511687157998545.csv"
$ws4.Range("B6").Value = "MM_stims-16511687158305962.csv"
$ws4.Range("B7").Value = "ZM_stims-16511687158154333.csv"

$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16511687159146957"
$ws5.Range("B2").Value = "vSAT_stims-1651168715867472.csv"
$ws5.Range("B3").Value = "SAT_stims-1651168715834595.csv"
$ws5.Range("B4").Value = "vSAT_stims-1651168715898432.csv"
$ws5.Range("B5").Value = "SAT_stims-16511687158514347.csv"
